$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 1
$ws.Range("B7").Value = 2
$ws.Range("C8").Value = 3
$ws.Range("D9").Value = 4
$ws.Range("E10").Value = 5

$ws.Range("E20").Select()
